$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "What is the smallest union territory in India?"
$ws.Range("B2").Value = "easy"
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 0

# Row 3
$ws.Range("A3").Value = "What does Lakshadweep mean in the local language?"
$ws.Range("F3").Value = 9
$ws.Range("G3").Value = -0.15625
$ws.Range("H3").Value = 0.34375

# Row 4
$ws.Range("A4").Value = "What is a permit to visit Lakshadweep in Kerala?"
$ws.Range("F4").Value = 10

# Row 5
$ws.Range("A5").Value = "How many of the Lakshadweep islands are inhabited by Muslims?"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 11

# Row 6
$ws.Range("A6").Value = "What are the main economic activities of a country?"
$ws.Range("E6").Value = 2

# Row 7
$ws.Range("A7").Value = "How many hectares of coconut are under cultivation?"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 9
$ws.Range("G7").Value = 0.5
$ws.Range("H7").Value = 0.5

$wb.Save()
